$wb = $excel.ActiveWorkbook

# --- SNMP sheet: securityName value khang -> public, selection move, column widths ---
$snmp = $wb.Worksheets.Item("SNMP")
$snmp.Range("A2").Value = "public"

$snmp.Columns.Item(2).ColumnWidth = 6.571428571428571
$snmp.Columns.Item(4).ColumnWidth = 5.857142857142857
$snmp.Columns.Item(6).ColumnWidth = 10.428571428571429

# --- New Alarm_Code sheet, inserted right after SNMP (becomes last sheet / active tab) ---
$ws = $wb.Worksheets.Add($null, $snmp)
$ws.Name = "Alarm_Code"

$ws.Range("A1").Value = "AlarmCode"
$ws.Range("B1").Value = "TrapID"
$ws.Range("D1").Value = "Message"
$ws.Range("A2").Value = "QADMN01005"
$ws.Range("C1").Value = "VerifyMessage"
$ws.Range("B2").Value = 17126
$ws.Range("C2").Value = "yes"
$ws.Range("D2").Value = "The information that you entered has been saved."

$ws.Columns.Item(1).ColumnWidth = 11.714285714285714
$ws.Columns.Item(2).ColumnWidth = 5.714285714285714
$ws.Columns.Item(3).ColumnWidth = 12.0
$ws.Columns.Item(4).ColumnWidth = 41.142857142857146

# Final selections: SNMP -> H2 (no longer active tab); Alarm_Code -> B2 (active tab)
$snmp.Range("H2").Select()
$ws.Select()
$ws.Range("B2").Select()
